$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (HeapSort): split "O(nlogn)" worst-case cell into a full
# best/average/worst set, adding the "all elements are equal" note ---
$ws.Range("D4").Value = "O(n)"
$ws.Range("E4").Value = "All elements are equal"
$ws.Range("F4").Value = "O(nlogn)"
$ws.Range("H4").Value = "O(nlogn)"

# --- Row 32 (Interval Scheduling): correct the recurrence + worst case ---
$ws.Range("B32").Value = "T(n) = T(n - 1) + O(nlogn)"
$ws.Range("F32").Value = "O(n^2)"
$ws.Range("H32").Value = "O(n^2)"

# --- New reviewer comments on the corrected worst-case cells ---
$noteText = "https://pdsaiitm.github.io/ says it's O(nlogn).   Need to confirm"
$ws.Range("F32").AddComment($noteText) | Out-Null
$ws.Range("H32").AddComment($noteText) | Out-Null

# --- Update the active selection left after editing ---
$ws.Range("B26").Select() | Out-Null
